# Auto-generated edit script applying crypto price/volume updates
# per commit: "Updated cryptos list on Mon Jun  3 04:15:06 UTC 2024 with GitHub Actions"
#
# D-column price values are stored as plain text in the workbook even though they
# look numeric (e.g. "1.00", "68.343.74" with thousands separators as dots). Excel
# auto-converts a numeric-looking string assigned via .Value into a real Number,
# which would silently drop things like trailing zeros ("1.00" -> 1). Prefixing
# with a leading apostrophe (Excel's classic 'treat as text' escape) keeps them as
# text, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.343.74"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "'3.800.01"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'608.28"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'163.32"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'3.798.46"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'6.96"
$ws.Range("E11").Value = "  +10.26%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "'35.08"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'4.436.76"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'3.820.19"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'68.358.36"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "'18.04"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'7.07"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'461.72"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'83.34"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'11.99"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'10.00"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -5.87%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").Value = "'28.96"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'9.06"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'0.146"
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("D39").Value = "'5.87"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "'0.975"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'3.17"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'43.22"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.297"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'153.03"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'46.84"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.40"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "'8.36"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "'380.18"
$ws.Range("E51").Value = "  -2.36%  "
